$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 6395
$wsExpo.Range("F4").Value = 5
$wsExpo.Range("F5").Value = 376
$wsExpo.Range("F8").Value = 530
$wsExpo.Range("F15").Value = 3112
$wsExpo.Range("F16").Value = 9
$wsExpo.Range("F18").Value = 1777

# Sheet "全部类型" (sheet4)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 6395
$wsAll.Range("F4").Value = 5
$wsAll.Range("F5").Value = 376
$wsAll.Range("F9").Value = 530
$wsAll.Range("F16").Value = 3112
$wsAll.Range("F17").Value = 9
$wsAll.Range("F19").Value = 1777
